$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "paiewise"

# "Sep-22-2023" would otherwise be auto-recognized as a date literal and
# converted to a serial date number. Build it as a text formula result in a
# scratch cell, then copy/paste-special the value only, which keeps it as
# literal text without stamping a new (Text) number-format style on B12.
$ws.Range("Z100").Formula = "=""Sep-22-2023"""
$ws.Range("Z100").Copy()
$ws.Range("B12").PasteSpecial(-4163)
$ws.Range("Z100").ClearContents()
$excel.CutCopyMode = $false

$ws.Range("C12").Value = "NV"
$ws.Range("D12").Value = "all_submitted_tracker_ninaSep-22-2023"

$ws.Range("D14").Select()
